$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Delete column V ("Nature") - this shifts every later column one to the left
#    (W->V, X->W, Y->X, Z->Y, AA->Z, AB->AA, AC->AB, AD->AC, AE->AD, AF->AE, AG->AF)
$ws.Columns.Item(22).Delete()

# 2. Apply the data-value edits (post-shift column letters)
$ws.Range("D2").Value = "Provide123*"
$ws.Range("D3").Value = "Provide123*"

$ws.Range("F2").Value = "Test Automation India Anmol"
$ws.Range("F3").Value = "Test Automation India Anmol"

$ws.Range("N2").Value = "Eur"

$ws.Range("E3").Value = "'010"
$ws.Range("G3").Value = "Product2"
$ws.Range("I3").Value = "HW - Network Security"
$ws.Range("J3").Value = "S00008044001"
$ws.Range("U3").Value = "C&S SUB CONTRACTOR"
$ws.Range("W3").Value = "'4201"
$ws.Range("X3").Value = "026909B2065"

$ws.Range("AD2").Value = "EMEAAD\srofidal"
$ws.Range("AE2").Value = "'"
$ws.Range("AD3").Value = "EMEAAD\srofidal"
$ws.Range("AE3").Value = "'"

# 3. Highlight the "PR tYPE" column (Q) with a yellow fill, matching the new styles added
$yellow = 65535
$ws.Range("Q1:Q3").Interior.Color = $yellow
$ws.Columns.Item(17).Interior.Color = $yellow

# 4. Cosmetic view updates matching the saved workbook window/selection state
$ws.Application.ActiveWindow.ScrollColumn = 19
$ws.Range("AF3").Select()
